# Scheduled runner update: refresh market price / profit data across Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 829.7
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 716.1667
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 716.1667
$ws.Range("M41").Value = -560
$ws.Range("N41").Value = -1596.1667

$ws.Range("H43").Value = 1160.909
$ws.Range("I43").Value = 840
$ws.Range("J43").Value = 1232.2222
$ws.Range("K43").Value = 840
$ws.Range("L43").Value = 1232.2222
$ws.Range("M43").Value = -771
$ws.Range("N43").Value = -1370.2222

$ws.Range("H51").Value = 2200
$ws.Range("I51").Value = 1800
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 1800
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -1316
$ws.Range("N51").Value = -3968

$ws.Range("H101").Value = 8578.223
$ws.Range("I101").Value = 683
$ws.Range("J101").Value = 12525.833
$ws.Range("K101").Value = 2049
$ws.Range("L101").Value = 37577.499
$ws.Range("M101").Value = -427
$ws.Range("N101").Value = -40821.499

$ws.Range("H113").Value = 2516
$ws.Range("I113").Value = 2525
$ws.Range("J113").Value = 2480
$ws.Range("K113").Value = 2525
$ws.Range("L113").Value = 2480
$ws.Range("M113").Value = 729
$ws.Range("N113").Value = -8988

$ws.Range("H129").Value = 898.4239
$ws.Range("I129").Value = 472.7
$ws.Range("J129").Value = 950.3415
$ws.Range("K129").Value = 1418.1
$ws.Range("L129").Value = 2851.0245
$ws.Range("M129").Value = 3581.9
$ws.Range("N129").Value = -12851.0245

$ws.Range("H132").Value = 1031.7435
$ws.Range("I132").Value = 1092.25
$ws.Range("J132").Value = 877.7273
$ws.Range("K132").Value = 3276.75
$ws.Range("L132").Value = 2633.1819
$ws.Range("M132").Value = -746.75
$ws.Range("N132").Value = -7693.1819

$ws.Range("H135").Value = 1621.8518
$ws.Range("I135").Value = 1720.6364
$ws.Range("J135").Value = 1187.2
$ws.Range("K135").Value = 15485.7276
$ws.Range("L135").Value = 10684.8
$ws.Range("M135").Value = -12950.7276
$ws.Range("N135").Value = -15754.8

$ws.Range("H136").Value = 27400
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 27400
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 27400
$ws.Range("N136").Value = -37600

$ws.Range("H137").Value = 1300.2572
$ws.Range("I137").Value = 1174
$ws.Range("J137").Value = 1575.7273
$ws.Range("K137").Value = 3522
$ws.Range("L137").Value = 4727.1819
$ws.Range("M137").Value = -972
$ws.Range("N137").Value = -9827.1819

$ws.Range("H138").Value = 4894.2173
$ws.Range("I138").Value = 912.12
$ws.Range("J138").Value = 9634.809999999999
$ws.Range("K138").Value = 2736.36
$ws.Range("L138").Value = 28904.43
$ws.Range("M138").Value = 2403.64
$ws.Range("N138").Value = -39184.43

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4243.4697
$ws.Range("I32").Value = 3043.6606
$ws.Range("J32").Value = 10962.4
$ws.Range("K32").Value = 3043.6606
$ws.Range("L32").Value = 10962.4
$ws.Range("M32").Value = -2756.6606
$ws.Range("N32").Value = -11536.4

$ws.Range("H45").Value = 13223.556
$ws.Range("I45").Value = 13626.5
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 13626.5
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -13249.5
$ws.Range("N45").Value = -10754

$ws.Range("H61").Value = 5981.2607
$ws.Range("I61").Value = 6709.55
$ws.Range("J61").Value = 1126
$ws.Range("K61").Value = 6709.55
$ws.Range("L61").Value = 1126
$ws.Range("M61").Value = -6497.55
$ws.Range("N61").Value = -1550

$ws.Range("H74").Value = 1591.95
$ws.Range("I74").Value = 1583.6875
$ws.Range("J74").Value = 1625
$ws.Range("K74").Value = 1583.6875
$ws.Range("L74").Value = 1625
$ws.Range("M74").Value = -709.6875
$ws.Range("N74").Value = -3373

$ws.Range("H77").Value = 1591.95
$ws.Range("I77").Value = 1583.6875
$ws.Range("J77").Value = 1625
$ws.Range("K77").Value = 7918.4375
$ws.Range("L77").Value = 8125
$ws.Range("M77").Value = -3550.4375
$ws.Range("N77").Value = -16861

$ws.Range("H122").Value = 25641024
$ws.Range("I122").Value = 25641024
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 76923072
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -76920622

$ws.Range("H132").Value = 3406.946
$ws.Range("I132").Value = 1924.8462
$ws.Range("J132").Value = 6910.091
$ws.Range("K132").Value = 5774.5386
$ws.Range("L132").Value = 20730.273
$ws.Range("M132").Value = -3244.5386
$ws.Range("N132").Value = -25790.273

$ws.Range("H136").Value = 5981.2607
$ws.Range("I136").Value = 6709.55
$ws.Range("J136").Value = 1126
$ws.Range("K136").Value = 20128.65
$ws.Range("L136").Value = 3378
$ws.Range("M136").Value = -17578.65
$ws.Range("N136").Value = -8478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 50780
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 50780
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 50780
$ws.Range("N57").Value = -52220

$ws.Range("H134").Value = 5313.9062
$ws.Range("I134").Value = 6474.773
$ws.Range("J134").Value = 2760
$ws.Range("K134").Value = 19424.319
$ws.Range("L134").Value = 8280
$ws.Range("M134").Value = -16889.319
$ws.Range("N134").Value = -13350

$ws.Range("H136").Value = 50780
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 50780
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 50780
$ws.Range("N136").Value = -60980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4056.1765
$ws.Range("I31").Value = 1493.9143
$ws.Range("J31").Value = 9661.125
$ws.Range("K31").Value = 1493.9143
$ws.Range("L31").Value = 9661.125
$ws.Range("M31").Value = -1198.9143
$ws.Range("N31").Value = -10251.125

$ws.Range("H34").Value = 4056.1765
$ws.Range("I34").Value = 1493.9143
$ws.Range("J34").Value = 9661.125
$ws.Range("K34").Value = 1493.9143
$ws.Range("L34").Value = 9661.125
$ws.Range("M34").Value = -1291.9143
$ws.Range("N34").Value = -10065.125

$ws.Range("H58").Value = 1437.95
$ws.Range("I58").Value = 1130.9231
$ws.Range("J58").Value = 2008.1428
$ws.Range("K58").Value = 1130.9231
$ws.Range("L58").Value = 2008.1428
$ws.Range("M58").Value = -927.9231
$ws.Range("N58").Value = -2414.1428

$ws.Range("H93").Value = 21500
$ws.Range("I93").Value = 8000
$ws.Range("J93").Value = 35000
$ws.Range("K93").Value = 8000
$ws.Range("L93").Value = 35000
$ws.Range("M93").Value = -6128
$ws.Range("N93").Value = -38744

$ws.Range("H132").Value = 2529.9722
$ws.Range("I132").Value = 2562.5
$ws.Range("J132").Value = 2416.125
$ws.Range("K132").Value = 7687.5
$ws.Range("L132").Value = 7248.375
$ws.Range("M132").Value = -5157.5
$ws.Range("N132").Value = -12308.375

$ws.Range("H134").Value = 4682.7856
$ws.Range("I134").Value = 5580.364
$ws.Range("J134").Value = 1391.6666
$ws.Range("K134").Value = 16741.092
$ws.Range("L134").Value = 4174.9998
$ws.Range("M134").Value = -14206.092
$ws.Range("N134").Value = -9244.9998

$ws.Range("H136").Value = 1437.95
$ws.Range("I136").Value = 1130.9231
$ws.Range("J136").Value = 2008.1428
$ws.Range("K136").Value = 3392.7693
$ws.Range("L136").Value = 6024.428400000001
$ws.Range("M136").Value = -842.7692999999999
$ws.Range("N136").Value = -11124.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3981.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3981.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 11944.5
$ws.Range("N55").Value = -12298.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2948109.8
$ws.Range("I122").Value = 3242720.5
$ws.Range("J122").Value = 2004
$ws.Range("K122").Value = 9728161.5
$ws.Range("L122").Value = 6012
$ws.Range("M122").Value = -9725711.5
$ws.Range("N122").Value = -10912

$ws.Range("H132").Value = 3935.353
$ws.Range("I132").Value = 5009.5
$ws.Range("J132").Value = 2980.5557
$ws.Range("K132").Value = 15028.5
$ws.Range("L132").Value = 8941.667099999999
$ws.Range("M132").Value = -12498.5
$ws.Range("N132").Value = -14001.6671

$ws.Range("H135").Value = 43991.668
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 43991.668
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 43991.668
$ws.Range("N135").Value = -54131.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1917191.6
$ws.Range("I22").Value = 5291478.5
$ws.Range("J22").Value = 2055.7568
$ws.Range("K22").Value = 5291478.5
$ws.Range("L22").Value = 2055.7568
$ws.Range("M22").Value = -5291183.5
$ws.Range("N22").Value = -2645.7568

$ws.Range("H27").Value = 1917191.6
$ws.Range("I27").Value = 5291478.5
$ws.Range("J27").Value = 2055.7568
$ws.Range("K27").Value = 5291478.5
$ws.Range("L27").Value = 2055.7568
$ws.Range("M27").Value = -5291371.5
$ws.Range("N27").Value = -2269.7568

$ws.Range("H132").Value = 13894458
$ws.Range("I132").Value = 19934100
$ws.Range("J132").Value = 3280.8
$ws.Range("K132").Value = 59802300
$ws.Range("L132").Value = 9842.400000000001
$ws.Range("M132").Value = -59799770
$ws.Range("N132").Value = -14902.4

$ws.Range("H136").Value = 8290.846
$ws.Range("I136").Value = 8337.625
$ws.Range("J136").Value = 8216
$ws.Range("K136").Value = 25012.875
$ws.Range("L136").Value = 24648
$ws.Range("M136").Value = -22462.875
$ws.Range("N136").Value = -29748

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1347.6
$ws.Range("I126").Value = 1347
$ws.Range("J126").Value = 1350
$ws.Range("K126").Value = 4041
$ws.Range("L126").Value = 4050
$ws.Range("M126").Value = -1571
$ws.Range("N126").Value = -8990

$ws.Range("H132").Value = 1453.0952
$ws.Range("I132").Value = 977.46155
$ws.Range("J132").Value = 2226
$ws.Range("K132").Value = 2932.38465
$ws.Range("L132").Value = 6678
$ws.Range("M132").Value = -402.38465
$ws.Range("N132").Value = -11738

$ws.Range("H136").Value = 3643.9092
$ws.Range("I136").Value = 4890.0835
$ws.Range("J136").Value = 2148.5
$ws.Range("K136").Value = 14670.2505
$ws.Range("L136").Value = 6445.5
$ws.Range("M136").Value = -12120.2505
$ws.Range("N136").Value = -11545.5

